# Generate Report for Handback
# Update the timestamp cells recorded in the handback status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 09:41:11"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 09:40:58"
$wsZhCn.Range("K2").Value = "2016-09-07 09:41:45"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-07 09:41:11"
$wsDeDe.Range("K2").Value = "2016-09-07 09:41:54"
